$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.099.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.911.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.97%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.57"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4848"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3825"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07366"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9354"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.83"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07816"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.931.73"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.517"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.641"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.65"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008843"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "28.143.89"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.88"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.165"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.170.44"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.28"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.911"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.57"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.125"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.32"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.955"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08949"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.341"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.254"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7752"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.688"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.629"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02060"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.104"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05322"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5499"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.985"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.034"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1528"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.476"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.67"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4835"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "107.82"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.72%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.660"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.44"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.50%  "